$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 1: headers (Testcases, Data1, Data2, Data3, Data4) ---
$ws1.Range("A1").Value = "Testcases"
$ws1.Range("B1").Value = "Data1"
$ws1.Range("C1").Value = "Data2"
$ws1.Range("D1").Value = "Data3"
$ws1.Range("E1").Value = "Data4"

# --- Row 2 ---
$ws1.Range("A2").Value = "Login test"
$ws1.Range("B2").Value = "2-2"
$ws1.Range("C2").Value = "2-3"
$ws1.Range("D2").Value = "2-4"
$ws1.Range("E2").Value = "2-5"

# --- Row 3 ---
$ws1.Range("A3").Value = "Purchase"
$ws1.Range("B3").Value = "3-2"
$ws1.Range("C3").Value = "3-3"
$ws1.Range("D3").Value = "3-4"
$ws1.Range("E3").Value = "3-5"

# --- Row 4 ---
$ws1.Range("A4").Value = "Add Profile"
$ws1.Range("B4").Value = "4-2"
$ws1.Range("C4").Value = "4-3"
$ws1.Range("D4").Value = "4-4"
$ws1.Range("E4").Value = "4-5"

# --- Row 5 ---
$ws1.Range("A5").Value = "Delete Profile"
$ws1.Range("B5").Value = "5-2"
$ws1.Range("C5").Value = "5-3"
$ws1.Range("D5").Value = "5-4"
$ws1.Range("E5").Value = "5-5"

# --- Row 6 ---
$ws1.Range("A6").Value = "Ignore Profile"
$ws1.Range("B6").Value = "6-2"
$ws1.Range("C6").Value = "6-3"
$ws1.Range("D6").Value = "6-4"
$ws1.Range("E6").Value = "6-5"

# Format the used range as Text (numFmtId 49 / "@"), matching the new style xf.
$ws1.Range("A1:E6").NumberFormat = "@"

# Selection state for sheet1
$ws1.Range("J7").Select()

# Selection state for sheet2 ("sample")
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("A1:E6").Select()

# Restore sheet1 as the active/visible tab
$ws1.Activate()

# Page setup (paper size + orientation) for sheet1
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
